$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits (rows 2-25, no row shifting involved yet) ---
$ws.Range("E5").Value = $null          # was -5
$ws.Range("E11").Value = -7.9          # was blank
$ws.Range("D19").Value = -15.5         # was blank
$ws.Range("E19").Value = $null         # was -6.5
$ws.Range("D21").Value = $null         # was -14.3
$ws.Range("D23").Value = -13.9         # was blank
$ws.Range("E25").Value = -7.1          # was blank

# --- Remove the "RM 232" row (original row 26) entirely ---
$ws.Rows.Item(26).Delete()

# --- Remove the "SC 92" row (now shifted up to row 27) entirely ---
$ws.Rows.Item(27).Delete()

# --- Remaining targeted value edits on rows that shifted up ---
# Row 27 is now "SC 101" (originally row 29); D27 goes blank
$ws.Range("D27").Value = $null         # was -14.6
# Row 29 is now "SC 119" (originally row 31); E29 goes blank
$ws.Range("E29").Value = $null         # was -6.8
# Row 33 is now "SC 232" (originally row 35); D33 gets a value
$ws.Range("D33").Value = -14.1         # was blank
